$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-126 down to 67-127.
$ws.Rows(66).Insert()

# Populate the newly inserted row 66 with the new weekly price record.
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44586
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112052
$ws.Range("G66").Value = "Albahaca"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 60
$ws.Range("K66").Value = 4000
$ws.Range("L66").Value = 4000
$ws.Range("M66").Value = 4000
$ws.Range("N66").Value = "$/docena de matas"
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 667
$ws.Range("Q66").Value = 6
$ws.Range("R66").Value = "Hortaliza"
